$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "FICHA" column: copy the header style from the last existing
# header (F1) into G1, then set the text, and set the data value in G2.
$ws.Range("F1").Copy($ws.Range("G1"))
$ws.Range("G1").Value = "FICHA"
$ws.Range("G2").Value = 2671143
$ws.Columns("G").ColumnWidth = 14.140625

# Update name in row 2 (drop surname "MOSQUERA")
$ws.Range("C2").Value = "MARLON"

# Remove the second data row (former row 3: 1108453116 / HP / 28)
$ws.Rows(3).Delete()

# Update selection to match target view state
$ws.Range("F3").Select()
